# Weekly update: insert a new price record for "Apio" (Femacal de La Calera)
# at row 289, pushing the existing rows 289-312 down to 290-313.
#
# The new record reuses most of the previous row's static attributes
# (market, region, product, quality, unit, origin, etc.) and carries new
# weekly figures for date / volume / min-max-avg prices / $ per kg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 289:312 down by inserting a fresh row at 289.
$ws.Rows.Item(289).Insert()

# Populate the newly inserted row 289 with this week's data.
$ws.Cells.Item(289, 1).Value  = 3
$ws.Cells.Item(289, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(289, 3).Value  = "Coquimbo"
$ws.Cells.Item(289, 4).Value  = 44578
$ws.Cells.Item(289, 5).Value  = 5
$ws.Cells.Item(289, 6).Value  = 100112017
$ws.Cells.Item(289, 7).Value  = "Apio"
$ws.Cells.Item(289, 8).Value  = "Americana (o)"
$ws.Cells.Item(289, 9).Value  = "Primera"
$ws.Cells.Item(289, 10).Value = 530
$ws.Cells.Item(289, 11).Value = 9000
$ws.Cells.Item(289, 12).Value = 9500
$ws.Cells.Item(289, 13).Value = 9236
$ws.Cells.Item(289, 14).Value = "$/docena de matas"
$ws.Cells.Item(289, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(289, 16).Value = 1539
$ws.Cells.Item(289, 17).Value = 6
$ws.Cells.Item(289, 18).Value = "Hortaliza"
